{"js": "// Insert three new \"List Paragraph\" bullet items (\"Pr\u00e1tica\",\n// \"Introdu\u00e7\u00e3o ao GitHub\", \"Push e Pull\") right after the existing\n// \"Comandos\" bullet item, matching the style/list numbering of the\n// surrounding items.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Comandos\" paragraph that the new items must follow.\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.trim() === \"Comandos\") {\n    anchor = p;\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find the \"Comandos\" paragraph.');\n}\n\n// Grab the list this paragraph belongs to so the new paragraphs can be\n// attached to the very same numbering definition (numId).\nconst list = anchor.list;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\nconst newTexts = [\"Pr\u00e1tica\", \"Introdu\u00e7\u00e3o ao GitHub\", \"Push e Pull\"];\n\nlet insertAfter = anchor;\nfor (const text of newTexts) {\n  const newPara = insertAfter.insertParagraph(text, \"After\");\n  newPara.style = \"List Paragraph\";\n  newPara.attachToList(listId, 0);\n  insertAfter = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Insert three new \"List Paragraph\" bullet items (\"Pr\u00e1tica\",\n# \"Introdu\u00e7\u00e3o ao GitHub\", \"Push e Pull\") right after the existing\n# \"Comandos\" bullet item, matching the style/list numbering of the\n# surrounding items.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByText {\n    param($doc, [string]$text)\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $t = $doc.Paragraphs.Item($i).Range.Text\n        # Paragraph.Range.Text carries the trailing paragraph/section mark\n        # (CR / cell mark / form-feed) - strip it before comparing.\n        $t = $t.TrimEnd([char]13, [char]7, [char]12)\n        if ($t -eq $text) { return $i }\n    }\n    return -1\n}\n\n$anchorIndex = Get-ParagraphIndexByText $d \"Comandos\"\nif ($anchorIndex -eq -1) {\n    throw 'Could not find the \"Comandos\" paragraph.'\n}\n\n# Grab the list that \"Comandos\" belongs to, so the new paragraphs can be\n# attached to that same numbering definition (numId).\n$listId = $d.Paragraphs.Item($anchorIndex).Range.ListFormat.List.ListID\n\n$newTexts = @(\"Pr\u00e1tica\", \"Introdu\u00e7\u00e3o ao GitHub\", \"Push e Pull\")\n\n$prevIndex = $anchorIndex\nforeach ($text in $newTexts) {\n    $prevParagraph = $d.Paragraphs.Item($prevIndex)\n    $prevParagraph.Range.InsertParagraphAfter()\n\n    $newIndex = $prevIndex + 1\n    $newParagraph = $d.Paragraphs.Item($newIndex)\n    # Attaches the same numbering (numId/ilvl) as the source list; the new\n    # paragraph already inherits the \"List Paragraph\" style from its\n    # predecessor, so no explicit Style= assignment is needed (and doing so\n    # would strip the numPr we just attached).\n    $newParagraph.AttachToList($listId, 0)\n    $newParagraph.Range.InsertAfter($text)\n\n    $prevIndex = $newIndex\n}\n"}
